$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Output $ws.Name
Write-Output $ws.Range("B5").Value2
Write-Output $ws.Range("C5").Value2
Write-Output $ws.Cells.Item(5,3).Value2
